$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.855.82'
$ws.Range('E2').Value = '  -0.15%  '
$ws.Range('D3').Value = '2.305.86'
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '102.60'
$ws.Range('E5').Value = '  +5.81%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '270.74'
$ws.Range('E6').Value = '  -0.14%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.626'
$ws.Range('E7').Value = '  -0.16%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.608'
$ws.Range('E9').Value = '  -2.61%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '45.90'
$ws.Range('E10').Value = '  +0.11%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0935'
$ws.Range('E11').Value = '  -1.74%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.94'
$ws.Range('E12').Value = '  -2.05%  '
$ws.Range('E13').Value = '  +1.67%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '15.88'
$ws.Range('E14').Value = '  +2.08%  '
$ws.Range('D15').Value = '2.651.65'
$ws.Range('E15').Value = '  -0.94%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.860'
$ws.Range('E16').Value = '  -0.83%  '
$ws.Range('D17').Value = '2.311.17'
$ws.Range('E17').Value = '  -0.82%  '
$ws.Range('D18').Value = '43.834.03'
$ws.Range('E18').Value = '  -0.23%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.0000110'
$ws.Range('E19').Value = '  +0.94%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.25'
$ws.Range('E20').Value = '  -2.68%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '72.37'
$ws.Range('E21').Value = '  -0.57%  '
$ws.Range('E22').Value = '  +8.34%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '233.62'
$ws.Range('E23').Value = '  -2.54%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.88'
$ws.Range('E24').Value = '  +14.01%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.19'
$ws.Range('E25').Value = '  -2.53%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.999'
$ws.Range('E26').Value = '  +0.00%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.23'
$ws.Range('E27').Value = '  -1.22%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '3.45'
$ws.Range('E28').Value = '  +0.35%  '
$ws.Range('E29').Value = '  +0.40%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '38.64'
$ws.Range('E30').Value = '  +0.69%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '177.44'
$ws.Range('E31').Value = '  +1.52%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '21.85'
$ws.Range('E32').Value = '  -2.52%  '
$ws.Range('E33').Value = '  -1.25%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.48'
$ws.Range('E34').Value = '  -0.09%  '
$ws.Range('E35').Value = '  -0.01%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.79'
$ws.Range('E36').Value = '  +8.32%  '
$ws.Range('E37').Value = '  -0.95%  '
$ws.Range('E38').Value = '  -2.09%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.59'
$ws.Range('E39').Value = '  +6.37%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.34'
$ws.Range('E40').Value = '  -0.87%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.235'
$ws.Range('E41').Value = '  -2.77%  '
$ws.Range('E42').Value = '  +0.61%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '12.23'
$ws.Range('E43').Value = '  +0.00%  '
$ws.Range('E44').Value = '  +4.11%  '
$ws.Range('E45').Value = '  -3.62%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '5.25'
$ws.Range('E46').Value = '  -2.66%  '
$ws.Range('E47').Value = '  -1.22%  '
$ws.Range('E48').Value = '  +0.91%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '98.50'
$ws.Range('E49').Value = '  -1.80%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.446'
$ws.Range('E50').Value = '  +7.22%  '
$ws.Range('E51').Value = '  +10.72%  '
